# "se realizan ajuste en tiempos" - refresh the test MSISDN/MSI numbers on
# sheet "Semilla 11" and add a new pending-test marker row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Semilla 11")

# --- Row 13 : brand-new pair (close to the old row-13 one) ---
$ws.Range("C13").Value2 = "3046010523"
$ws.Range("D13").Value2 = "732111193280544"

# --- Row 12 : brand-new pair ---
$ws.Range("C12").Value2 = "3052749177"
$ws.Range("D12").Value2 = "732111324709512"

# --- Row 10 : brand-new MSISDN/MSI pair + new "Cedula Cliente" ---
$ws.Range("B10").Value2 = "484303795"
$ws.Range("C10").Value2 = "3046010569"
$ws.Range("D10").Value2 = "732111193280551"
$ws.Range("E10").Value2 = "3046008593"

# --- Row 9 : reuse the (now-freed) numbers previously on row 12 ---
$ws.Range("C9").Value2 = "3045987650"
$ws.Range("D9").Value2 = "732111193278858"

# --- Row 11 : same pair as row 9 ---
$ws.Range("C11").Value2 = "3045987650"
$ws.Range("D11").Value2 = "732111193278858"

# --- Row 14 : reuse old row-9 "E" value, mark it visually (underline) ---
$ws.Range("C14").Value2 = "3045984556"
$ws.Range("C14").Font.Underline = -4119

# --- New row 17 : empty placeholder cell with a bordered style ---
$ws.Range("B17").Borders.Item(1).LineStyle = -4142

# --- View bookkeeping: scroll back to top-left and land the selection on B17 ---
$ws.Range("A1").Select()
$ws.Range("B17").Select()
